$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Sprachangaben")
$src.Copy($null, $src)
$ws = $wb.Worksheets.Item(6)
$ws.Rows(13).Copy($ws.Rows(14))
$ws.Rows(13).Copy($ws.Rows(15))
Write-Host "done"
